$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell H1 ("S") ---
$ws.Range("H1").Value = "S"

# --- Row 2: measured values updated, formulas introduced ---
$ws.Range("A2").Value = 220
$ws.Range("B2").Formula = "=A2"
$ws.Range("D2").Value = 210
$ws.Range("E2").Formula = "=D2/1000"
$ws.Range("F2").Formula = "=2*C2/(E2^2)"
$ws.Range("H2").Formula = "=(9.8*E2^2)/2"
$ws.Range("I2").Formula = "=H2*1000"

# --- Row 3 ---
$ws.Range("D3").Value = 285
$ws.Range("B3").Formula = "=B2+A3"
$ws.Range("E3:E5").Formula = "=D3/1000"
$ws.Range("F3").Formula = "=2*C3/(E3^2)"
$ws.Range("H3:H5").Formula = "=(9.8*E3^2)/2"
$ws.Range("I3:I5").Formula = "=H3*1000"

# touch J3:K5 so the (empty, but formatted) cells are persisted
$ws.Range("J3:K5").Borders.LineStyle = -4142

# --- Row 4 ---
$ws.Range("D4").Value = 352
$ws.Range("B4:B5").Formula = "=B3+A4"
$ws.Range("F4:F5").Formula = "=2*C4/(E4^2)"

# --- Row 5 ---
$ws.Range("D5").Value = 403

# --- Row 7: average ---
$ws.Range("E7").Value = "Média"
$ws.Range("F2").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("F7").Formula = "=AVERAGE(F2:F5)"

# --- Selection matches the diff ---
$ws.Range("J2:K7").Select()

# --- Page setup (printable area defaults) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
